$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.351.69"
$ws.Range("E2").Value = "  -2.44%  "

$ws.Range("D3").Value = "'3.390.79"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'574.43"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "'136.12"
$ws.Range("E6").Value = "  +7.82%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "'3.390.06"
$ws.Range("E8").Value = "  -0.54%  "

$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("D10").Value = "'7.60"
$ws.Range("E10").Value = "  +2.73%  "

$ws.Range("E11").Value = "  +2.03%  "

$ws.Range("E12").Value = "  +2.50%  "

$ws.Range("D13").Value = "'3.964.36"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").Value = "'0.0000175"
$ws.Range("E15").Value = "  +0.89%  "

$ws.Range("D16").Value = "'3.383.81"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("D17").Value = "'25.31"
$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "'61.439.02"
$ws.Range("E18").Value = "  -2.39%  "

$ws.Range("D19").Value = "'14.08"
$ws.Range("E19").Value = "  +6.84%  "

$ws.Range("E20").Value = "  +2.31%  "

$ws.Range("D21").Value = "'9.46"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("D22").Value = "'376.94"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").Value = "'0.569"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("D24").Value = "'3.521.86"
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.12%  "

$ws.Range("D26").Value = "'70.71"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("E27").Value = "  +9.36%  "

$ws.Range("D28").Value = "'1.70"
$ws.Range("E28").Value = "  +21.73%  "

$ws.Range("D29").Value = "'7.75"
$ws.Range("E29").Value = "  +10.72%  "

$ws.Range("D30").Value = "'0.988"
$ws.Range("E30").Value = "  -1.02%  "

$ws.Range("D31").Value = "'8.17"
$ws.Range("E31").Value = "  +3.81%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.156"
$ws.Range("E32").Value = "  +3.43%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "'2.16"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").Value = "'3.421.67"
$ws.Range("E35").Value = "  -0.55%  "

$ws.Range("D36").Value = "'23.53"
$ws.Range("E36").Value = "  +3.21%  "

$ws.Range("D37").Value = "'5.60"
$ws.Range("E37").Value = "  +5.77%  "

$ws.Range("E38").Value = "  +6.02%  "

$ws.Range("D39").Value = "'6.96"
$ws.Range("E39").Value = "  +3.40%  "

$ws.Range("D40").Value = "'162.54"
$ws.Range("E40").Value = "  -1.93%  "

$ws.Range("D41").Value = "'0.0790"
$ws.Range("E41").Value = "  +4.06%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("E43").Value = "  +12.74%  "

$ws.Range("D44").Value = "'4.44"
$ws.Range("E44").Value = "  +3.48%  "

$ws.Range("D45").Value = "'41.50"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "'0.762"
$ws.Range("E46").Value = "  -1.85%  "

$ws.Range("D47").Value = "'1.62"
$ws.Range("E47").Value = "  +2.89%  "

$ws.Range("D48").Value = "'23.72"
$ws.Range("E48").Value = "  +3.49%  "

$ws.Range("D49").Value = "'6.97"
$ws.Range("E49").Value = "  +4.39%  "

$ws.Range("E50").Value = "  +12.67%  "

$ws.Range("D51").Value = "'0.902"
$ws.Range("E51").Value = "  +5.11%  "
